# Insert a new weekly record for "Orégano" (Mercado Mayorista Lo Valledor de
# Santiago) above the existing most-recent entry. This pushes the existing
# rows 321:369 down to 322:370 and fills the freed-up row 321 with the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 321:369 down to 322:370, leaving row 321 blank (ready for the
# new record) and growing the used range to A1:R370.
$ws.Rows("321:321").Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Range("A321").Value = 6
$ws.Range("B321").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C321").Value = "Metropolitana"
$ws.Range("D321").Value = 45218
$ws.Range("E321").Value = 13
$ws.Range("F321").Value = 100112029
$ws.Range("G321").Value = "Orégano"
$ws.Range("H321").Value = "Sin especificar"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 30
$ws.Range("K321").Value = 16000
$ws.Range("L321").Value = 16000
$ws.Range("M321").Value = 16000
$ws.Range("N321").Value = "$/docena de atados"
$ws.Range("O321").Value = "Región Metropolitana"
$ws.Range("P321").Value = 5333
$ws.Range("Q321").Value = 3
$ws.Range("R321").Value = "Hortaliza"
